$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 6 ("SOURCES OF FINANCE") gets a new built-in table style.
#    (a:tableStyleId {1F2CE20C-...} -> {A8FF3916-...})
# ---------------------------------------------------------------------------
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $targetSlide = $slide
            $targetShape = $shape
        }
    }
}

if ($targetShape -ne $null) {
    $table = $targetShape.Table
    $table.ApplyStyle("{A8FF3916-FE92-4139-9BC5-DDE24FCE578E}")
}

# ---------------------------------------------------------------------------
# 2) Swap the two theme colour palettes: the slide master's theme (the
#    "Integral" greens/yellows) is replaced with the plain default "Office"
#    palette, while the other theme part keeps the Integral colours.
#    Only the RGB slots are reachable through the exposed object model, so
#    drive the swap through SlideMaster.Theme.ThemeColorScheme.Colors(i).RGB.
# ---------------------------------------------------------------------------
function ToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeHex = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ToRGB($officeHex[$i - 1])
}
